$d = $word.ActiveDocument

$oldText = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Locate the astromap credit/link sentence, spanning several differently
# formatted runs (plain text, "(", hyperlink-styled URL, ")."), and
# replace it as a whole so the whole sentence collapses into a single run
# with the updated (2022) year in the URL.
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute($oldText)

if ($found) {
    $target = $find.Parent
    $target.Delete()
    $target.InsertAfter($newText)
}

$found
